$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix "player" -> "Player" (A34) per commit message
$ws.Range('A34').Value = 'Player'

# Populate new SPANISH translations column (C)
$ws.Range('C1').Value = 'SPANISH'
$ws.Range('C2').Value = 'Ninguno'
$ws.Range('C3').Value = 'Suelo'
$ws.Range('C4').Value = 'Suelo Cyan'
$ws.Range('C5').Value = 'Suelo Rojo'
$ws.Range('C6').Value = 'Suelo Naranja'
$ws.Range('C7').Value = 'Suelo Grande'
$ws.Range('C8').Value = 'Suelo 2'
$ws.Range('C9').Value = 'Pared'
$ws.Range('C10').Value = 'Pared Sin Color'
$ws.Range('C11').Value = 'Pared X'
$ws.Range('C12').Value = 'Ventana'
$ws.Range('C13').Value = 'Luz Direccional'
$ws.Range('C14').Value = 'Luz Puntual'
$ws.Range('C15').Value = 'Lámpara'
$ws.Range('C16').Value = 'Ventilación Verde'
$ws.Range('C17').Value = 'Ventilación Cyan'
$ws.Range('C18').Value = 'Paquete de Salud'
$ws.Range('C19').Value = 'Paquete de Munición'
$ws.Range('C20').Value = 'Sierra'
$ws.Range('C21').Value = 'Waypoint de Sierra'
$ws.Range('C22').Value = 'Interruptor'
$ws.Range('C23').Value = 'Spawn del Jugador'
$ws.Range('C24').Value = 'Cubo'
$ws.Range('C25').Value = 'Laser'
$ws.Range('C26').Value = 'Trampa de Fuego'
$ws.Range('C27').Value = 'Colisión'
$ws.Range('C28').Value = 'Trigger Final'
$ws.Range('C29').Value = 'Placa de Presión'
$ws.Range('C30').Value = 'Pantalla'
$ws.Range('C31').Value = 'Pantalla Pequeña'
$ws.Range('C32').Value = 'Ventana Frágil'
$ws.Range('C33').Value = 'Trigger'
$ws.Range('C34').Value = 'Jugador'
$ws.Range('C36').Value = 'Cuando Se Activa'
$ws.Range('C37').Value = 'Cuando Se Desactiva'
$ws.Range('C38').Value = 'Cuando Se Invierte'

# Set column C width to match column in target (23 chars)
$ws.Columns.Item(3).ColumnWidth = 22 + 1/6

# Give C1 header the same bold/yellow/bordered/centered style as A1:B1
$ws.Range('A1').Copy()
$ws.Range('C1').PasteSpecial(-4122)
$ws.Range('C1').Value = 'SPANISH'

# Rows 36-38 (When Activating/Deactivating/Inverting) lose their special
# bold/yellow-highlight style and become plain cells like the rest of the sheet
$ws.Range('A2').Copy()
$ws.Range('A36:B38').PasteSpecial(-4122)
$ws.Range('A36').Value = 'WhenActivatingEvents'
$ws.Range('B36').Value = 'When Activating'
$ws.Range('A37').Value = 'WhenDeactivatingEvents'
$ws.Range('B37').Value = 'When Deactivating'
$ws.Range('A38').Value = 'WhenInvertingEvents'
$ws.Range('B38').Value = 'When Inverting'

$excel.CutCopyMode = $false

# Put selection on D38 to mirror final cursor position after data entry
$ws.Range('D38').Select()
